$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Internet Retail "
$ws.Range("B7").Value = 4

$ws.Range("A8").Value = "Biotechnology "
$ws.Range("B8").Value = 3

$ws.Range("A9").Value = "Telecom Services "
$ws.Range("B9").Value = 2

$ws.Range("A10").Value = "Agricultural Inputs "
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "Oil Gas E&amp;P "
$ws.Range("B11").Value = 1

$ws.Range("A12").Value = "Exchange Traded Fund "
$ws.Range("B12").Value = 1

$ws.Range("A13").Value = "Lodging "
$ws.Range("B13").Value = 1
